$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text changes: capitalize / rename the "location", "service",
#     and "practitioner" headers (N1/O1/P1) to Practitioner/Location/HealthCareService.
#     Setting these three values removes the now-unreferenced "location"/"service"
#     shared strings and appends the new ones, which is exactly what the target
#     workbook's sharedStrings table reflects (practitioner data cells are untouched,
#     only the header text itself changes).
$ws.Range("N1").Value = "Practitioner"
$ws.Range("O1").Value = "Location"
$ws.Range("P1").Value = "HealthCareService"

# --- Column width changes: new col P (width ~16.16 chars) and col Q widened to 18.5
$ws.Columns.Item(16).ColumnWidth = 15.333333333333334
$ws.Columns.Item(17).ColumnWidth = 17.666666666666668

# --- Date number format: admit_date / discharge_date columns switch from the
#     built-in datetime format (numFmtId 22) to the built-in short-date format
#     (numFmtId 14, canonical format code "mm-dd-yy"). Apply to one cell, then
#     propagate via a formats-only paste so every cell in the range shares the
#     same style entry instead of each getting its own duplicate.
$ws.Range("Q2").NumberFormat = "mm-dd-yy"
$ws.Range("Q2").Copy() | Out-Null
$ws.Range("Q2:R11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Selection / scroll position
$ws.Range("I16").Select() | Out-Null
